# Add a new "Flumazenil" row to the DOSES table, inserted right before the
# existing "Fluoxetine" row (i.e. becomes the new row 28, pushing every
# subsequent row down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 28 (shifts Fluoxetine ... Carbamazepine down by one).
$ws.Rows("28:28").Insert()

# Fill in the new row's three columns.
$ws.Range("A28").Value = 'Flumazenil'
$ws.Range("B28").Value = '"Varies by indication, administered IV. **Reversal of Sedation (Adult):** Initial dose 0.2 mg over 15 seconds. May repeat 0.2 mg at 60-second intervals up to a maximum total dose of 1 mg. **Suspected Overdose (Adult):** Initial dose 0.2 mg over 30 seconds. If no response after 30 seconds, administer 0.3 mg. If no response, subsequent doses of 0.5 mg at 1-minute intervals up to a maximum total cumulative dose of 3 mg."'
$ws.Range("C28").Value = '"Benzodiazepine Antagonist. Used to reverse sedation from benzodiazepines. Its short half-life (40-80 min) often leads to **resedation**, requiring repeat dosing or a continuous IV infusion. Use is **contraindicated** in patients with suspected serious tricyclic antidepressant overdose, or in those relying on benzodiazepines to control a life-threatening condition (e.g., status epilepticus) due to the risk of precipitating seizures."'

# Match the saved selection state (cell C28 selected).
$ws.Range("C28").Select() | Out-Null
